# ATPT_YR_FIN.xlsx - "Doing Updates for Financials"
# Insert a new first data-column (D) in front of the existing period columns
# on the ATPT sheet, shifting the old D:K columns to E:L, then populate the
# brand-new column D with the latest reporting period (29-Jun-2018) and the
# handful of line items that have real figures for that period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATPT")

# --- Insert a new column before D, shifting D:K -> E:L -------------------
$ws.Columns("D").Insert()

# Match the column width Excel carried over to the (now shifted) E column so
# the freshly inserted D column looks the same as its neighbours.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# --- "Period Ending" header rows: new latest period = 29-Jun-2018 --------
# (only rows 7/38/80 actually carry the header dates; set format + value
# together so we never touch rows that have no D:K cells at all.)
$headerRows = @(7,38,80)
foreach ($r in $headerRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "[$-409]d\-mmm\-yy;@"
    $cell.Value2 = 43281
}

# --- Data rows: default the new period's figures to 0 ---------------------
# Only the rows that actually have the D:K data block get a new value here;
# pure section-title rows (e.g. 36/37/78/79) have no D:K cells at all and
# must stay untouched.
$defaultRows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,`
                 39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,61,62,63,64,65,`
                 67,68,69,70,71,72,73,74,75,77,`
                 81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102)

foreach ($r in $defaultRows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.NumberFormat = "#,##0"
    $cell.Value2 = 0
}

# --- Data rows with real figures for the new period -----------------------
$ws.Range("D59").Value2 = 100     # Other Current Liabilities
$ws.Range("D60").Value2 = 100     # Total Current Liabilities
$ws.Range("D66").Value2 = 100     # Total Liabilities
$ws.Range("D76").Value2 = -100    # Total Stockholder Equity
